# Apply rail car trace report update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test_format_trace")

# Row 1 - search description text (date/time updated)
$ws.Range("A1").Value = "Description unknown, completed 06/21/2023 09:18:52 EDT, by WPJTOWN1.The search returned: 8 events."

# Row 3
$ws.Range("A3").Value = "MWCX"
$ws.Range("B3").Value = 100705
$ws.Range("C3").Value = "AMORY"
$ws.Range("D3").Value = "MS"
$ws.Range("F3").Value = 20
$ws.Range("G3").Value = 1759
$ws.Range("H3").Value = "Arrive In-Transit"
$ws.Range("I3").Value = "RHLA85"
$ws.Range("J3").Value = "LOVELAND"
$ws.Range("K3").Value = "CO"
$ws.Range("L3").Value = 267061
$ws.Range("M3").Value = 72200
$ws.Range("N3").Value = 194861
$ws.Range("O3").Value = "MWCX100705"

# Row 4
$ws.Range("A4").Value = "MWCX"
$ws.Range("B4").Value = 100715
$ws.Range("C4").Value = "ARKWRIGHT"
$ws.Range("D4").Value = "GA"
$ws.Range("F4").Value = 21
$ws.Range("G4").Value = 603
$ws.Range("H4").Value = "Departure"
$ws.Range("J4").Value = "LOVELAND"
$ws.Range("K4").Value = "CO"
$ws.Range("L4").Value = 267600
$ws.Range("M4").Value = 73600
$ws.Range("N4").Value = 194000
$ws.Range("O4").Value = "MWCX100715"

# Row 5
$ws.Range("A5").Value = "MWCX"
$ws.Range("B5").Value = 102553
$ws.Range("C5").Value = "COLORADO SPRIN"
$ws.Range("D5").Value = "CO"
$ws.Range("F5").Value = 21
$ws.Range("G5").Value = 714
$ws.Range("H5").Value = "Departure"
$ws.Range("I5").Value = "HKCKDE"
$ws.Range("J5").Value = "LOVELAND"
$ws.Range("K5").Value = "CO"
$ws.Range("L5").Value = 281050
$ws.Range("M5").Value = 73400
$ws.Range("N5").Value = 207650
$ws.Range("O5").Value = "MWCX102553"

# Row 6
$ws.Range("A6").Value = "ITFX"
$ws.Range("B6").Value = 9728
$ws.Range("C6").Value = "JOHNSTOWN"
$ws.Range("D6").Value = "CO"
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1812
$ws.Range("H6").Value = "Placed Actual"
$ws.Range("J6").Value = "JOHNSTOWN"
$ws.Range("K6").Value = "CO"
$ws.Range("L6").Value = 202950
$ws.Range("M6").Value = 0
$ws.Range("N6").Value = 202950
$ws.Range("O6").Value = "ITFX9728"

# Row 7
$ws.Range("A7").Value = "MWCX"
$ws.Range("B7").Value = 102276
$ws.Range("C7").Value = "JOHNSTOWN"
$ws.Range("D7").Value = "CO"
$ws.Range("F7").Value = 12
$ws.Range("G7").Value = 1304
$ws.Range("H7").Value = "Placed Actual"
$ws.Range("J7").Value = "LOVELAND"
$ws.Range("K7").Value = "CO"
$ws.Range("L7").Value = 280350
$ws.Range("M7").Value = 78900
$ws.Range("N7").Value = 201450
$ws.Range("O7").Value = "MWCX102276"

# Row 8
$ws.Range("A8").Value = "MWCX"
$ws.Range("B8").Value = 102166
$ws.Range("C8").Value = "JOHNSTOWN"
$ws.Range("D8").Value = "CO"
$ws.Range("F8").Value = 12
$ws.Range("G8").Value = 1304
$ws.Range("H8").Value = "Placed Actual"
$ws.Range("I8").Value = ""
$ws.Range("J8").Value = "LOVELAND"
$ws.Range("K8").Value = "CO"
$ws.Range("L8").Value = 282400
$ws.Range("M8").Value = 82000
$ws.Range("N8").Value = 200400
$ws.Range("O8").Value = "MWCX102166"

# Row 9
$ws.Range("A9").Value = "MWCX"
$ws.Range("B9").Value = 102330
$ws.Range("C9").Value = "JOHNSTOWN"
$ws.Range("D9").Value = "CO"
$ws.Range("F9").Value = 15
$ws.Range("G9").Value = 1435
$ws.Range("H9").Value = "Placed Actual"
$ws.Range("I9").Value = ""
$ws.Range("J9").Value = "LOVELAND"
$ws.Range("K9").Value = "CO"
$ws.Range("L9").Value = 284850
$ws.Range("M9").Value = 79300
$ws.Range("N9").Value = 205550
$ws.Range("O9").Value = "MWCX102330"

# Row 10
$ws.Range("A10").Value = "MWCX"
$ws.Range("B10").Value = 102328
$ws.Range("C10").Value = "MEMPHIS"
$ws.Range("D10").Value = "TN"
$ws.Range("F10").Value = 20
$ws.Range("G10").Value = 1950
$ws.Range("H10").Value = "Bad Order"
$ws.Range("I10").Value = "L 000"
$ws.Range("J10").Value = "LOVELAND"
$ws.Range("K10").Value = "CO"
$ws.Range("L10").Value = 280550
$ws.Range("M10").Value = 79500
$ws.Range("N10").Value = 201050
$ws.Range("O10").Value = "MWCX102328"
